# "Sorts formulae by depth"
#
# Calculation sheet: the lookup table used by the C2 VLOOKUP was shrunk
# from K2:L4 to K2:L2 (really K2:L3, 2 rows), and the "No"/default row
# that used to live at row 4 (K4/L4/M4) was moved up to row 3
# (K3/L3/M3). The B4:B13 helper formulas, which used the workbook-level
# defined name NEEDS, now reference Calculation!$B$3 directly, so the
# NEEDS defined name is no longer needed and gets removed.

$wb = $excel.ActiveWorkbook

# --- Workbook level: drop the now-unused "NEEDS" defined name ---------
foreach ($n in @($wb.Names)) {
    if ($n.Name -eq "NEEDS") {
        $n.Delete()
    }
}

# --- Calculation sheet --------------------------------------------------
$calc = $wb.Worksheets.Item("Calculation")

# Lookup range shrunk from K2:L4 to K2:L3
$calc.Range("C2").Formula = "=VLOOKUP(B2,K2:L3,2,FALSE)"

# The "no" row moves from row 4 up to row 3
$calc.Range("K3").Value = 0
$calc.Range("L3").Value = 0
$calc.Range("M3").Value = "no"
$calc.Range("K4:M4").ClearContents()

# B4:B13 formulas now reference Calculation!$B$3 directly instead of
# going through the NEEDS defined name
$digits = @("0", "1", "2", "3", "4", "5", "6", "7", "8", "9")
for ($i = 0; $i -lt $digits.Length; $i++) {
    $row = 4 + $i
    $digit = $digits[$i]
    $calc.Cells.Item($row, 2).Formula = "=IF(IFERROR(SEARCH(""$digit"", Calculation!`$B`$3), 0),1,0)"
}

# C5:C13 keep the same per-row formula
for ($row = 5; $row -le 13; $row++) {
    $calc.Cells.Item($row, 3).Formula = "=IF(B$row>0,10,0)"
}

# Update the sheet's active selection/scroll position
$calc.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 2
$calc.Range("C2").Select()
